$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("hits")

# --- New row block 1: rows 41-49 (meta_0..meta_3, rint_0..rint_3 for AX-11150762) ---
# Row 41
$ws.Range("A41").Value = "STUDY"
$ws.Range("B41").Value = "SNP"
$ws.Range("C41").Value = "CHR"
$ws.Range("D41").Value = "BP"
$ws.Range("E41").Value = "MAF"
$ws.Range("F41").Value = "A1"
$ws.Range("G41").Value = "A2"
$ws.Range("H41").Value = "DIR"
$ws.Range("I41").Value = "BETA"
$ws.Range("J41").Value = "SE_BETA"
$ws.Range("K41").Value = "Q"
$ws.Range("L41").Value = "P"
# Row 42
$ws.Range("A42").Value = "meta_0"
$ws.Range("B42").Value = "AX-11150762"
$ws.Range("C42").Value = 1
$ws.Range("D42").Value = 55496039
$ws.Range("E42").Value = 0.05637
$ws.Range("F42").Value = "T"
$ws.Range("G42").Value = "C"
$ws.Range("H42").Value = "+++-+  ++"
$ws.Range("I42").Value = 0.0428
$ws.Range("J42").Value = 0.0139
$ws.Range("K42").Value = 0.03899
$ws.Range("L42").Value = 0.002149
# Row 43
$ws.Range("A43").Value = "meta_1"
$ws.Range("B43").Value = "AX-11150762"
$ws.Range("C43").Value = 1
$ws.Range("D43").Value = 55496039
$ws.Range("E43").Value = 0.05637
$ws.Range("F43").Value = "T"
$ws.Range("G43").Value = "C"
$ws.Range("H43").Value = "+++-+  ++"
$ws.Range("I43").Value = 0.047
$ws.Range("J43").Value = 0.0139
$ws.Range("K43").Value = 0.03806
$ws.Range("L43").Value = 0.0006894
# Row 44
$ws.Range("A44").Value = "meta_2"
$ws.Range("B44").Value = "AX-11150762"
$ws.Range("C44").Value = 1
$ws.Range("D44").Value = 55496039
$ws.Range("E44").Value = 0.05637
$ws.Range("F44").Value = "T"
$ws.Range("G44").Value = "C"
$ws.Range("H44").Value = "+++-+  ++"
$ws.Range("I44").Value = 0.0411
$ws.Range("J44").Value = 0.0139
$ws.Range("K44").Value = 0.03258
$ws.Range("L44").Value = 0.003251
# Row 45
$ws.Range("A45").Value = "meta_3"
$ws.Range("B45").Value = "AX-11150762"
$ws.Range("C45").Value = 1
$ws.Range("D45").Value = 55496039
$ws.Range("E45").Value = 0.05637
$ws.Range("F45").Value = "T"
$ws.Range("G45").Value = "C"
$ws.Range("H45").Value = "+++-+  ++"
$ws.Range("I45").Value = 0.037
$ws.Range("J45").Value = 0.014
$ws.Range("K45").Value = 0.02969
$ws.Range("L45").Value = 0.007931
# Row 46
$ws.Range("A46").Value = "rint_0"
$ws.Range("B46").Value = "AX-11150762"
$ws.Range("C46").Value = 1
$ws.Range("D46").Value = 55496039
$ws.Range("E46").Value = 0.05637
$ws.Range("F46").Value = "T"
$ws.Range("G46").Value = "C"
$ws.Range("H46").Value = "NA"
$ws.Range("I46").Value = 0.06388
$ws.Range("J46").Value = 0.02052
$ws.Range("K46").Value = "NA"
$ws.Range("L46").Value = 0.001854
# Row 47
$ws.Range("A47").Value = "rint_1"
$ws.Range("B47").Value = "AX-11150762"
$ws.Range("C47").Value = 1
$ws.Range("D47").Value = 55496039
$ws.Range("E47").Value = 0.05637
$ws.Range("F47").Value = "T"
$ws.Range("G47").Value = "C"
$ws.Range("H47").Value = "NA"
$ws.Range("I47").Value = 0.07066
$ws.Range("J47").Value = 0.02043
$ws.Range("K47").Value = "NA"
$ws.Range("L47").Value = 0.0005426
# Row 48
$ws.Range("A48").Value = "rint_2"
$ws.Range("B48").Value = "AX-11150762"
$ws.Range("C48").Value = 1
$ws.Range("D48").Value = 55496039
$ws.Range("E48").Value = 0.05637
$ws.Range("F48").Value = "T"
$ws.Range("G48").Value = "C"
$ws.Range("H48").Value = "NA"
$ws.Range("I48").Value = 0.03684
$ws.Range("J48").Value = 0.02138
$ws.Range("K48").Value = "NA"
$ws.Range("L48").Value = 0.08488
# Row 49
$ws.Range("A49").Value = "rint_3"
$ws.Range("B49").Value = "AX-11150762"
$ws.Range("C49").Value = 1
$ws.Range("D49").Value = 55496039
$ws.Range("E49").Value = 0.05637
$ws.Range("F49").Value = "T"
$ws.Range("G49").Value = "C"
$ws.Range("H49").Value = "NA"
$ws.Range("I49").Value = 0.03921
$ws.Range("J49").Value = 0.02145
$ws.Range("K49").Value = "NA"
$ws.Range("L49").Value = 0.06753

# --- New row block 2: rows 51-59 (meta_0..meta_3, rint_0..rint_3 for AX-39911995) ---
# Row 51
$ws.Range("A51").Value = "STUDY"
$ws.Range("B51").Value = "SNP"
$ws.Range("C51").Value = "CHR"
$ws.Range("D51").Value = "BP"
$ws.Range("E51").Value = "MAF"
$ws.Range("F51").Value = "A1"
$ws.Range("G51").Value = "A2"
$ws.Range("H51").Value = "DIR"
$ws.Range("I51").Value = "BETA"
$ws.Range("J51").Value = "SE_BETA"
$ws.Range("K51").Value = "Q"
$ws.Range("L51").Value = "P"
# Row 52
$ws.Range("A52").Value = "meta_0"
$ws.Range("B52").Value = "AX-39911995"
$ws.Range("C52").Value = 1
$ws.Range("D52").Value = 55504650
$ws.Range("E52").Value = 0.3181
$ws.Range("F52").Value = "A"
$ws.Range("G52").Value = "G"
$ws.Range("H52").Value = "---+-  --"
$ws.Range("I52").Value = -0.0064
$ws.Range("J52").Value = 0.0069
$ws.Range("K52").Value = 0.7365
$ws.Range("L52").Value = 0.353
# Row 53
$ws.Range("A53").Value = "meta_1"
$ws.Range("B53").Value = "AX-39911995"
$ws.Range("C53").Value = 1
$ws.Range("D53").Value = 55504650
$ws.Range("E53").Value = 0.3181
$ws.Range("F53").Value = "A"
$ws.Range("G53").Value = "G"
$ws.Range("H53").Value = "---+-  --"
$ws.Range("I53").Value = -0.0148
$ws.Range("J53").Value = 0.0069
$ws.Range("K53").Value = 0.7518
$ws.Range("L53").Value = 0.03265
# Row 54
$ws.Range("A54").Value = "meta_2"
$ws.Range("B54").Value = "AX-39911995"
$ws.Range("C54").Value = 1
$ws.Range("D54").Value = 55504650
$ws.Range("E54").Value = 0.3181
$ws.Range("F54").Value = "A"
$ws.Range("G54").Value = "G"
$ws.Range("H54").Value = "---+-  --"
$ws.Range("I54").Value = -0.0119
$ws.Range("J54").Value = 0.007
$ws.Range("K54").Value = 0.868
$ws.Range("L54").Value = 0.08902
# Row 55
$ws.Range("A55").Value = "meta_3"
$ws.Range("B55").Value = "AX-39911995"
$ws.Range("C55").Value = 1
$ws.Range("D55").Value = 55504650
$ws.Range("E55").Value = 0.3181
$ws.Range("F55").Value = "A"
$ws.Range("G55").Value = "G"
$ws.Range("H55").Value = "---+-  --"
$ws.Range("I55").Value = -0.0142
$ws.Range("J55").Value = 0.007
$ws.Range("K55").Value = 0.8512
$ws.Range("L55").Value = 0.042
# Row 56
$ws.Range("A56").Value = "rint_0"
$ws.Range("B56").Value = "AX-39911995"
$ws.Range("C56").Value = 1
$ws.Range("D56").Value = 55504650
$ws.Range("E56").Value = 0.3181
$ws.Range("F56").Value = "A"
$ws.Range("G56").Value = "G"
$ws.Range("H56").Value = "NA"
$ws.Range("I56").Value = -0.009374
$ws.Range("J56").Value = 0.01023
$ws.Range("K56").Value = "NA"
$ws.Range("L56").Value = 0.3596
# Row 57
$ws.Range("A57").Value = "rint_1"
$ws.Range("B57").Value = "AX-39911995"
$ws.Range("C57").Value = 1
$ws.Range("D57").Value = 55504650
$ws.Range("E57").Value = 0.3181
$ws.Range("F57").Value = "A"
$ws.Range("G57").Value = "G"
$ws.Range("H57").Value = "NA"
$ws.Range("I57").Value = -0.02243
$ws.Range("J57").Value = 0.01022
$ws.Range("K57").Value = "NA"
$ws.Range("L57").Value = 0.02817
# Row 58
$ws.Range("A58").Value = "rint_2"
$ws.Range("B58").Value = "AX-39911995"
$ws.Range("C58").Value = 1
$ws.Range("D58").Value = 55504650
$ws.Range("E58").Value = 0.3181
$ws.Range("F58").Value = "A"
$ws.Range("G58").Value = "G"
$ws.Range("H58").Value = "NA"
$ws.Range("I58").Value = -0.007956
$ws.Range("J58").Value = 0.01053
$ws.Range("K58").Value = "NA"
$ws.Range("L58").Value = 0.45
# Row 59
$ws.Range("A59").Value = "rint_3"
$ws.Range("B59").Value = "AX-39911995"
$ws.Range("C59").Value = 1
$ws.Range("D59").Value = 55504650
$ws.Range("E59").Value = 0.3181
$ws.Range("F59").Value = "A"
$ws.Range("G59").Value = "G"
$ws.Range("H59").Value = "NA"
$ws.Range("I59").Value = -0.009969
$ws.Range("J59").Value = 0.01057
$ws.Range("K59").Value = "NA"
$ws.Range("L59").Value = 0.3456

# --- Column widths (best-fit sizing for the new table columns) ---
$ws.Columns.Item(1).ColumnWidth = 6.333333333333334
$ws.Columns.Item(2).ColumnWidth = 11.0
$ws.Columns.Item(3).ColumnWidth = 4.0
$ws.Columns.Item(4).ColumnWidth = 8.166666666666668
$ws.Columns.Item(5).ColumnWidth = 7.166666666666666
$ws.Columns.Item(6).ColumnWidth = 2.5
$ws.Columns.Item(7).ColumnWidth = 2.5
$ws.Columns.Item(8).ColumnWidth = 8.666666666666668
$ws.Columns.Item(9).ColumnWidth = 7.166666666666666
$ws.Columns.Item(10).ColumnWidth = 8.666666666666668
$ws.Columns.Item(11).ColumnWidth = 7.166666666666666
$ws.Columns.Item(12).ColumnWidth = 9.166666666666668

# --- Restore the active selection to match the edited area ---
$ws.Range("B51").Select() | Out-Null
